# "Support of Dione for Ubuntu" — add a duplicate "Ubuntu 24.04" column (H)
# to the Raspberry Pi sheet, mirroring column F (Ubuntu 24.04 on Rpi4) but
# flagging the non-Dione camera modules as "Not supported" on this new
# target, and make the Raspberry Pi sheet the active one with H5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)       # "Raspberry Pi"
$ws2 = $wb.Worksheets.Item(2)      # "Flex cable"

# --- widen the merged header cell from G1:G2 to G1:H2 ------------------
$ws.Range("G1:H2").Merge()

# --- restore the header formatting (Merge() above recomputes styles) ---
$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- new column H data rows 3-6: duplicate column F's formatting -------
$ws.Range("F3:F6").Copy()
$ws.Range("H3:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- new column H data rows 7-9: duplicate column G's formatting -------
# (these rows render the "Not supported" text like column G, not the
#  "OK" style used in column F)
$ws.Range("G7:G9").Copy()
$ws.Range("H7:H9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- values for the new column ------------------------------------------
$ws.Range("H3").Value2 = "Ubuntu 24.04"
$ws.Range("H4").Value2 = "6.8.0-1004"
$ws.Range("H5").Value2 = "OK"
$ws.Range("H6").Value2 = "OK"
$ws.Range("H7").Value2 = "Not supported"
$ws.Range("H8").Value2 = "Not supported"
$ws.Range("H9").Value2 = "Not supported"

# --- row heights tightened slightly for the three module rows ----------
$ws.Rows.Item(7).RowHeight = 18
$ws.Rows.Item(8).RowHeight = 19.5
$ws.Rows.Item(9).RowHeight = 18

# --- give the new column a sensible width -------------------------------
$ws.Columns.Item(8).ColumnWidth = 12.5

# --- make "Raspberry Pi" the active sheet / tab, select H5 --------------
$ws.Activate()
$ws.Range("H5").Select()

# --- the "Flex cable" sheet is no longer the active tab -----------------
$ws2.Range("A9").Select()
$ws.Activate()
